$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Five bullet points describing an exam|ple from the NYFD of:"
#    -> "Five bullet points describing an example from the NYFD of:"
#    The bookmark that currently sits between "exam" and "ple" is
#    subsumed by this edit (it disappears from this paragraph).
# ---------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute(
    "Five bullet points describing an example from the NYFD of:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Five bullet points describing an example from the NYFD of:", 2)

# ---------------------------------------------------------------------
# 2) "...sought by O'Hagan and Lindsay would"
#    -> "...sought by John O'Hagan and John Lindsay would"
# ---------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute(
    "sought by O’Hagan and Lindsay would",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "sought by John O’Hagan and John Lindsay would", 2)

# ---------------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark at the new edit location: right
#    before "Lindsay" (i.e. between the inserted "John " and "Lindsay
#    would"), matching where Word leaves it after the last keystroke.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$r3 = $d.Content
$found3 = $r3.Find.Execute(
    "Lindsay would", $true, $false, $false, $false, $false, $true, 1,
    $false, "", 0)

$gobackRange = $d.Range($r3.Start, $r3.Start)
$d.Bookmarks.Add("_GoBack", $gobackRange)

Write-Output "found1=$found1 found2=$found2 found3=$found3"
